# The deck currently applies the "Integral" theme (ppt/theme/theme2.xml) to
# the slide master, while the unused "Office Theme" colors sit in
# ppt/theme/theme1.xml (only wired to the Notes Master). The edit swaps the
# two themes' content, so the slide master (and therefore every slide) ends
# up styled with the standard Office Theme color scheme, Arial-based font
# scheme already matches between the two themes so only the color scheme
# needs to move.
#
# PowerPoint's legacy ColorScheme object (8 well-known slots plus the extra
# accent5/accent6/hyperlink/followed-hyperlink slots reachable by index)
# writes straight into the active theme part, so we drive the swap through
# $p.SlideMaster.ColorScheme.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
